# Add "Всього" (Total) and percentage columns (L, M) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): L1 = "Всього", M1 = 0.05 ---
$ws.Range("L1").Value = "Всього"
$ws.Range("M1").Value = 0.05

# Match the look of the other header cells (K1 uses the green "total" style)
# by copying its formatting onto the two new header cells.
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (2-27): L = weighted total ("Всього"), M = that total's percentage share ---
$totals = @{
    2  = @(747,   37.35)
    3  = @(882,   44.1)
    4  = @(3132,  156.6)
    5  = @(3195,  159.75)
    6  = @(4185,  209.25)
    7  = @(4455,  222.75)
    8  = @(5562,  278.1)
    9  = @(5688,  284.4)
    10 = @(5985,  299.25)
    11 = @(7731,  386.55)
    12 = @(11124, 556.2)
    13 = @(11349, 567.45)
    14 = @(11889, 594.45)
    15 = @(12879, 643.95)
    16 = @(13410, 670.5)
    17 = @(14130, 706.5)
    18 = @(14238, 711.9)
    19 = @(14274, 713.7)
    20 = @(14355, 717.75)
    21 = @(14409, 720.45)
    22 = @(15021, 751.05)
    23 = @(15066, 753.3)
    24 = @(15291, 764.55)
    25 = @(15336, 766.8)
    26 = @(15444, 772.2)
    27 = @(15507, 775.35)
}

foreach ($r in 2..27) {
    $pair = $totals[$r]
    $ws.Cells.Item($r, 12).Value = $pair[0]    # column L
    $ws.Cells.Item($r, 13).Value = $pair[1]    # column M
}
